$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 10499.5
$ws.Range("I34").Value = 10499.5
$ws.Range("K34").Value = 10499.5
$ws.Range("M34").Value = -10296.5
$ws.Range("H36").Value = 10499.5
$ws.Range("I36").Value = 10499.5
$ws.Range("K36").Value = 10499.5
$ws.Range("M36").Value = -9784.5
$ws.Range("H62").Value = 1950.5
$ws.Range("I62").Value = 2048.25
$ws.Range("J62").Value = 1755
$ws.Range("K62").Value = 2048.25
$ws.Range("L62").Value = 1755
$ws.Range("M62").Value = -1424.25
$ws.Range("N62").Value = -3003
$ws.Range("H65").Value = 1950.5
$ws.Range("I65").Value = 2048.25
$ws.Range("J65").Value = 1755
$ws.Range("K65").Value = 10241.25
$ws.Range("L65").Value = 8775
$ws.Range("M65").Value = -7121.25
$ws.Range("N65").Value = -15015
$ws.Range("H98").Value = 1073.1428
$ws.Range("I98").Value = 832.7917
$ws.Range("J98").Value = 1597.5454
$ws.Range("K98").Value = 832.7917
$ws.Range("L98").Value = 1597.5454
$ws.Range("M98").Value = 665.2083
$ws.Range("N98").Value = -4593.5454
$ws.Range("H105").Value = 48999.5
$ws.Range("J105").Value = 48999.5
$ws.Range("L105").Value = 48999.5
$ws.Range("N105").Value = -55987.5
$ws.Range("H109").Value = 68000
$ws.Range("J109").Value = 68000
$ws.Range("L109").Value = 68000
$ws.Range("N109").Value = -70774
$ws.Range("H122").Value = 1073.1428
$ws.Range("I122").Value = 832.7917
$ws.Range("J122").Value = 1597.5454
$ws.Range("K122").Value = 2498.3751
$ws.Range("L122").Value = 4792.6362
$ws.Range("M122").Value = -48.3751000000002
$ws.Range("N122").Value = -9692.636200000001
$ws.Range("H135").Value = 55556708
$ws.Range("I135").Value = 1142.3334
$ws.Range("J135").Value = 166667840
$ws.Range("K135").Value = 10281.0006
$ws.Range("L135").Value = 1500010560
$ws.Range("M135").Value = -7746.000599999999
$ws.Range("N135").Value = -1500015630
$ws.Range("H137").Value = 1557.9032
$ws.Range("I137").Value = 1330.3043
$ws.Range("J137").Value = 2212.25
$ws.Range("K137").Value = 3990.9129
$ws.Range("L137").Value = 6636.75
$ws.Range("M137").Value = -1440.9129
$ws.Range("N137").Value = -11736.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4816.6416
$ws.Range("I32").Value = 3408.5652
$ws.Range("K32").Value = 3408.5652
$ws.Range("M32").Value = -3121.5652
$ws.Range("H61").Value = 1741.238
$ws.Range("I61").Value = 1621.2
$ws.Range("J61").Value = 2041.3334
$ws.Range("K61").Value = 1621.2
$ws.Range("L61").Value = 2041.3334
$ws.Range("M61").Value = -1409.2
$ws.Range("N61").Value = -2465.3334
$ws.Range("H74").Value = 2367.2
$ws.Range("I74").Value = 993
$ws.Range("K74").Value = 993
$ws.Range("M74").Value = -119
$ws.Range("H77").Value = 2367.2
$ws.Range("I77").Value = 993
$ws.Range("K77").Value = 4965
$ws.Range("M77").Value = -597
$ws.Range("H111").Value = 43000
$ws.Range("J111").Value = 43000
$ws.Range("L111").Value = 43000
$ws.Range("N111").Value = -51180
$ws.Range("H136").Value = 1741.238
$ws.Range("I136").Value = 1621.2
$ws.Range("J136").Value = 2041.3334
$ws.Range("K136").Value = 4863.6
$ws.Range("L136").Value = 6124.0002
$ws.Range("M136").Value = -2313.6
$ws.Range("N136").Value = -11224.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5557.2
$ws.Range("I134").Value = 6767.737
$ws.Range("J134").Value = 1723.8334
$ws.Range("K134").Value = 20303.211
$ws.Range("L134").Value = 5171.5002
$ws.Range("M134").Value = -17768.211
$ws.Range("N134").Value = -10241.5002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2506.3572
$ws.Range("I31").Value = 2711.625
$ws.Range("K31").Value = 2711.625
$ws.Range("M31").Value = -2416.625
$ws.Range("H34").Value = 2506.3572
$ws.Range("I34").Value = 2711.625
$ws.Range("K34").Value = 2711.625
$ws.Range("M34").Value = -2509.625
$ws.Range("H105").Value = 1089.0834
$ws.Range("I105").Value = 1105.625
$ws.Range("K105").Value = 1105.625
$ws.Range("M105").Value = 641.375
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -52620
$ws.Range("H132").Value = 1852.4
$ws.Range("I132").Value = 1467.3846
$ws.Range("K132").Value = 4402.1538
$ws.Range("M132").Value = -1872.1538
$ws.Range("H134").Value = 1803.4615
$ws.Range("I134").Value = 1703.75
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 5111.25
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -2576.25
$ws.Range("N134").Value = -14070
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1428771.2
$ws.Range("I4").Value = 1666799.9
$ws.Range("K4").Value = 5000399.699999999
$ws.Range("M4").Value = -5000287.699999999
$ws.Range("H23").Value = 153.88889
$ws.Range("I23").Value = 181
$ws.Range("J23").Value = 99.666664
$ws.Range("K23").Value = 543
$ws.Range("L23").Value = 298.999992
$ws.Range("M23").Value = -308
$ws.Range("N23").Value = -768.999992
$ws.Range("H131").Value = 18560.902
$ws.Range("J131").Value = 19473.514
$ws.Range("L131").Value = 58420.542
$ws.Range("N131").Value = -68500.542
$ws.Range("H132").Value = 1506.8
$ws.Range("I132").Value = 1165
$ws.Range("J132").Value = 1734.6666
$ws.Range("K132").Value = 10485
$ws.Range("L132").Value = 15611.9994
$ws.Range("M132").Value = -7955
$ws.Range("N132").Value = -20671.9994
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 25505000
$ws.Range("I10").Value = 25505000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 25505000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -25504831
$ws.Range("N10").ClearContents()
$ws.Range("H113").Value = 1240.7778
$ws.Range("J113").Value = 1302.1666
$ws.Range("L113").Value = 1302.1666
$ws.Range("N113").Value = -5642.1666
$ws.Range("H126").Value = 2177867
$ws.Range("I126").Value = 11114733
$ws.Range("J126").Value = 50041.668
$ws.Range("K126").Value = 33344199
$ws.Range("L126").Value = 150125.004
$ws.Range("M126").Value = -33341729
$ws.Range("N126").Value = -155065.004
$ws.Range("H127").Value = 37220.5
$ws.Range("J127").Value = 37220.5
$ws.Range("L127").Value = 37220.5
$ws.Range("N127").Value = -47140.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2364.611
$ws.Range("I7").Value = 2100.0833
$ws.Range("J7").Value = 2893.6667
$ws.Range("K7").Value = 2100.0833
$ws.Range("L7").Value = 2893.6667
$ws.Range("M7").Value = -1988.0833
$ws.Range("N7").Value = -3117.6667
$ws.Range("H61").Value = 2280.7058
$ws.Range("I61").Value = 1952
$ws.Range("J61").Value = 2883.3333
$ws.Range("K61").Value = 1952
$ws.Range("L61").Value = 2883.3333
$ws.Range("M61").Value = -1750
$ws.Range("N61").Value = -3287.3333
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 2280.7058
$ws.Range("I113").Value = 1952
$ws.Range("J113").Value = 2883.3333
$ws.Range("K113").Value = 1952
$ws.Range("L113").Value = 2883.3333
$ws.Range("M113").Value = 218
$ws.Range("N113").Value = -7223.3333
$ws.Range("H126").Value = 2364.611
$ws.Range("I126").Value = 2100.0833
$ws.Range("J126").Value = 2893.6667
$ws.Range("K126").Value = 6300.249899999999
$ws.Range("L126").Value = 8681.000100000001
$ws.Range("M126").Value = -3830.249899999999
$ws.Range("N126").Value = -13621.0001
$ws.Range("H132").Value = 2072.7693
$ws.Range("I132").Value = 1300.5
$ws.Range("J132").Value = 2304.45
$ws.Range("K132").Value = 3901.5
$ws.Range("L132").Value = 6913.349999999999
$ws.Range("M132").Value = -1371.5
$ws.Range("N132").Value = -11973.35
$ws.Range("H136").Value = 5120.1
$ws.Range("I136").Value = 4101.1875
$ws.Range("K136").Value = 12303.5625
$ws.Range("M136").Value = -9753.5625
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2437.2
$ws.Range("I126").Value = 2013.4166
$ws.Range("J126").Value = 4132.3335
$ws.Range("K126").Value = 6040.2498
$ws.Range("L126").Value = 12397.0005
$ws.Range("M126").Value = -3570.2498
$ws.Range("N126").Value = -17337.0005
$ws.Range("H135").Value = 136547.5
$ws.Range("J135").Value = 136547.5
$ws.Range("L135").Value = 136547.5
$ws.Range("N135").Value = -146687.5
$ws.Range("H136").Value = 17363274
$ws.Range("I136").Value = 32682050
$ws.Range("K136").Value = 98046150
$ws.Range("M136").Value = -98043600
